$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.052.85'
$ws.Range('E2').Value = '  -4.51%  '
$ws.Range('D3').Value = '2.619.06'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '516.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.51%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('E10').Value = '  -3.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.334'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = '3.077.11'
$ws.Range('E13').Value = '  -3.20%  '
$ws.Range('D14').Value = '58.046.24'
$ws.Range('E14').Value = '  -4.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').Value = '2.618.32'
$ws.Range('E17').Value = '  -10.51%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '334.14'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.49%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.12%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '63.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.422'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -3.12%  '
$ws.Range('E26').Value = '  +0.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.77%  '
$ws.Range('D28').Value = '0.0₃0778'
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.23%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.33%  '
$ws.Range('B35').Value = 'SuiNetwork'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.899'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.16'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.841'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('E39').Value = '  -6.40%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  -1.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0962'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '267.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.69%  '
$ws.Range('E45').Value = '  +1.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.08%  '
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('D48').Value = '2.024.02'
$ws.Range('E48').Value = '  -5.57%  '
$ws.Range('E49').Value = '  -3.09%  '
$ws.Range('E50').Value = '  -5.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.33%  '
